$d = $word.ActiveDocument

# 1. Append the new sentence to the Description/Overview paragraph.
$rng = $d.Content
$found = $rng.Find.Execute("information is displayed in the boxes of information that is editable and can be save or unchanged", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "could not find target sentence" }
$rng.Collapse(0)
$rng.InsertAfter(". There is a selection box with the list of robots to be edited and text boxes to edit the robots name, team, and code. 4 buttons are available to save changes, cancel changes, create a new robot, or go back to the main menu.")

# 2. Delete the old _GoBack bookmark (Word's "last edit location" marker),
#    which used to sit right after the constructor write-up paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
}

# 3. Re-add _GoBack right after the text we just inserted (collapsed / zero length),
#    working around a paragraph-end Range resolution quirk by anchoring on a
#    temporary placeholder character that is removed afterwards.
$rng.Collapse(0)
$rng.InsertAfter("X")
$placeholderStart = $rng.Start
$bmRange = $d.Range($placeholderStart, $placeholderStart)
$d.Bookmarks.Add("_GoBack", $bmRange)
$delRange = $d.Range($placeholderStart, $placeholderStart + 1)
$delRange.Text = ""
